{"js": "// Replace the division-problem answers in the table with the new set.\n// Each old text is unique within the document, so a plain literal search\n// (matchCase) + full-range replace is unambiguous.\nconst replacements = [\n  [\"780\u00f78=97, 4\", \"106\u00f74=26, 2\"],\n  [\"541\u00f78=67, 5\", \"266\u00f76=44, 2\"],\n  [\"694\u00f72=347, 0\", \"940\u00f77=134, 2\"],\n  [\"166\u00f76=27, 4\", \"647\u00f79=71, 8\"],\n  [\"437\u00f77=62, 3\", \"871\u00f73=290, 1\"],\n  [\"327\u00f75=65, 2\", \"830\u00f72=415, 0\"],\n  [\"203\u00f77=29, 0\", \"304\u00f79=33, 7\"],\n  [\"792\u00f74=198, 0\", \"837\u00f74=209, 1\"],\n  [\"266\u00f78=33, 2\", \"460\u00f75=92, 0\"],\n  [\"369\u00f76=61, 3\", \"958\u00f78=119, 6\"],\n  [\"568\u00f78=71, 0\", \"664\u00f72=332, 0\"],\n  [\"261\u00f75=52, 1\", \"146\u00f75=29, 1\"],\n  [\"125\u00f79=13, 8\", \"901\u00f79=100, 1\"],\n  [\"734\u00f72=367, 0\", \"681\u00f76=113, 3\"],\n  [\"435\u00f76=72, 3\", \"545\u00f72=272, 1\"],\n  [\"463\u00f72=231, 1\", \"569\u00f72=284, 1\"],\n  [\"512\u00f78=64, 0\", \"745\u00f75=149, 0\"],\n  [\"985\u00f75=197, 0\", \"809\u00f75=161, 4\"],\n  [\"296\u00f78=37, 0\", \"148\u00f79=16, 4\"],\n  [\"648\u00f73=216, 0\", \"280\u00f79=31, 1\"],\n  [\"951\u00f74=237, 3\", \"414\u00f79=46, 0\"],\n  [\"547\u00f73=182, 1\", \"901\u00f76=150, 1\"],\n  [\"904\u00f75=180, 4\", \"317\u00f73=105, 2\"],\n  [\"961\u00f78=120, 1\", \"573\u00f76=95, 3\"],\n  [\"737\u00f72=368, 1\", \"579\u00f77=82, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"780\u00f78=97, 4\";    new = \"106\u00f74=26, 2\"},\n    @{old = \"541\u00f78=67, 5\";    new = \"266\u00f76=44, 2\"},\n    @{old = \"694\u00f72=347, 0\";   new = \"940\u00f77=134, 2\"},\n    @{old = \"166\u00f76=27, 4\";    new = \"647\u00f79=71, 8\"},\n    @{old = \"437\u00f77=62, 3\";    new = \"871\u00f73=290, 1\"},\n    @{old = \"327\u00f75=65, 2\";    new = \"830\u00f72=415, 0\"},\n    @{old = \"203\u00f77=29, 0\";    new = \"304\u00f79=33, 7\"},\n    @{old = \"792\u00f74=198, 0\";   new = \"837\u00f74=209, 1\"},\n    @{old = \"266\u00f78=33, 2\";    new = \"460\u00f75=92, 0\"},\n    @{old = \"369\u00f76=61, 3\";    new = \"958\u00f78=119, 6\"},\n    @{old = \"568\u00f78=71, 0\";    new = \"664\u00f72=332, 0\"},\n    @{old = \"261\u00f75=52, 1\";    new = \"146\u00f75=29, 1\"},\n    @{old = \"125\u00f79=13, 8\";    new = \"901\u00f79=100, 1\"},\n    @{old = \"734\u00f72=367, 0\";   new = \"681\u00f76=113, 3\"},\n    @{old = \"435\u00f76=72, 3\";    new = \"545\u00f72=272, 1\"},\n    @{old = \"463\u00f72=231, 1\";   new = \"569\u00f72=284, 1\"},\n    @{old = \"512\u00f78=64, 0\";    new = \"745\u00f75=149, 0\"},\n    @{old = \"985\u00f75=197, 0\";   new = \"809\u00f75=161, 4\"},\n    @{old = \"296\u00f78=37, 0\";    new = \"148\u00f79=16, 4\"},\n    @{old = \"648\u00f73=216, 0\";   new = \"280\u00f79=31, 1\"},\n    @{old = \"951\u00f74=237, 3\";   new = \"414\u00f79=46, 0\"},\n    @{old = \"547\u00f73=182, 1\";   new = \"901\u00f76=150, 1\"},\n    @{old = \"904\u00f75=180, 4\";   new = \"317\u00f73=105, 2\"},\n    @{old = \"961\u00f78=120, 1\";   new = \"573\u00f76=95, 3\"},\n    @{old = \"737\u00f72=368, 1\";   new = \"579\u00f77=82, 5\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
